$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column E (old E = "fantasy points" shifts to G)
$ws.Range("E1:F1").EntireColumn.Insert()

# New header cells
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Fill height/weight values for the 16 data rows (rows 2-17)
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 255
}
